$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.031.62"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +2.95%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.633.13"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +2.38%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.00%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'595.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.50%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'156.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +4.99%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.02%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +1.08%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +8.10%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.400"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +5.01%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +0.66%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +2.37%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'29.13"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +6.45%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  +21.69%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'3.110.83"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +2.58%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'64.979.94"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +3.04%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.635.24"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.03%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'12.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +3.44%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'4.80"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'352.72"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +2.48%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'7.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +8.77%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  +0.20%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'68.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +2.18%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +0.22%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'9.52"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +5.00%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -0.31%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.165"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.94%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +1.26%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'0.0₃0948"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +11.35%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -0.13%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("B31").Value = "'Bittensor"
$ws.Range("B31").Style = "Normal"
$ws.Range("C31").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("C31").Style = "Normal"
$ws.Range("D31").Value = "'523.39"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -5.66%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("B32").Value = "'PancakeSwap"
$ws.Range("B32").Style = "Normal"
$ws.Range("C32").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").Value = "'2.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +4.39%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'1.77"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +2.10%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'5.59"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +8.34%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'6.32"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +5.92%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +3.10%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'20.26"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +4.73%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'163.36"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -1.08%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +6.25%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +0.19%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -0.05%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'42.26"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +6.62%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'165.38"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.09%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'4.08"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +3.58%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.0613"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +4.69%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  +2.52%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  +8.40%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +3.31%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  +3.22%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0978"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +2.03%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'19.44"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +3.07%  "
$ws.Range("E51").Style = "Normal"
